$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$workpackage = "EU Exit PPP Regulation Service"
$user = "m1006990"
$caseId = "Case 0"
$success = "Success"
$notes = "Test Case"

$msgUpload119 = "Could not find the UI element corresponding to this selector:" + [char]10 + "<webctrl id='id__119' tag='SPAN'/>" + [char]10 + "The closest matches found are:" + [char]10 + "[93%] <webctrl id='id__159' tag='SPAN'/>" + [char]10 + "[86%] <webctrl id='id__142' tag='SPAN'/>" + [char]10 + "[86%] <webctrl id='id__147' tag='SPAN'/>" + [char]10 + "[86%] <webctrl id='id__150' tag='SPAN'/>" + [char]10 + "[86%] <webctrl id='id__153' tag='SPAN'/>" + [char]10 + "[86%] <webctrl id='id__156' tag='SPAN'/>" + [char]10 + "[86%] <webctrl id='id__186' tag='SPAN'/>" + [char]10 + "[83%] <webctrl id='id__0' tag='SPAN'/>" + [char]10 + "[64%] <webctrl id='header58-modifiedByColumn_1102' tag='SPAN'/>" + [char]10 + "[62%] <webctrl id='header58-modifiedByColumn_1102-name' tag='SPAN'/> at Source: Invoke Upload File - Amended workflow: Click 'Upload'"

$msgOutsideBounds = "Cannot send input to UI element because it is outside of screen bounds. at Source: Invoke Upload File - Amended workflow: Click 'Upload'"

$crlf = [char]13 + [char]10
$msgMoveFile = "Could not find the user-interface (UI) element for this action." + $crlf + $crlf + "Possible solutions:" + $crlf + " " + [char]8226 + "  Ensure application is opened and the UI element is visible on the screen at execution time" + $crlf + " " + [char]8226 + "  Edit the Target of the UI activity and use Validation to debug the issue. " + $crlf + " " + [char]8226 + "  If needed, re-indicate the element as its properties might have changed" + $crlf + " " + [char]8226 + "  Use ""Check state"" activity to check the application state before executing the action" + $crlf + " " + [char]8226 + "  Increase the ""Delay before"" value to allow time to the application to render entirely and become responsive at Source: Invoke Move File - Amended workflow: Click Document to Move"

$msgUpload504 = "Could not find the UI element corresponding to this selector:" + [char]10 + "<webctrl id='id__504' tag='SPAN'/>" + [char]10 + "The closest matches found are:" + [char]10 + "[93%] <webctrl id='id__150' tag='SPAN'/>" + [char]10 + "[92%] <webctrl id='id__0' tag='SPAN'/>" + [char]10 + "[86%] <webctrl id='id__142' tag='SPAN'/>" + [char]10 + "[86%] <webctrl id='id__147' tag='SPAN'/>" + [char]10 + "[86%] <webctrl id='id__153' tag='SPAN'/>" + [char]10 + "[86%] <webctrl id='id__156' tag='SPAN'/>" + [char]10 + "[86%] <webctrl id='id__159' tag='SPAN'/>" + [char]10 + "[79%] <webctrl id='id__186' tag='SPAN'/>" + [char]10 + "[64%] <webctrl id='header58-displayNameColumn_504' tag='SPAN'/>" + [char]10 + "[63%] <webctrl id='header58-dateModifiedColumn_506' tag='SPAN'/> at Source: Invoke Upload File - Amended workflow: Click 'Upload'"

$rows = @(
    @{ Row=2; Start=44615.466597222221; End=44615.466597222221; Total=44615.466597222221; Msg=$msgUpload119;   Wrap=$true },
    @{ Row=3; Start=44615.470902777779; End=44615.470902777779; Total=44615.470902777779; Msg=$msgOutsideBounds; Wrap=$false },
    @{ Row=4; Start=44615.493020833332; End=44615.493020833332; Total=44615.493009259262; Msg=$msgMoveFile;    Wrap=$true },
    @{ Row=5; Start=44615.501631944448; End=44615.501631944448; Total=44615.501620370371; Msg=$msgOutsideBounds; Wrap=$false },
    @{ Row=6; Start=44615.503993055558; End=44615.503993055558; Total=44615.503993055558; Msg=$msgOutsideBounds; Wrap=$false },
    @{ Row=7; Start=44615.510810185187; End=44615.510810185187; Total=44615.510798611111; Msg=$msgOutsideBounds; Wrap=$false },
    @{ Row=8; Start=44615.519965277781; End=44615.519965277781; Total=44615.519953703704; Msg=$msgUpload504;   Wrap=$true },
    @{ Row=9; Start=44615.522592592592; End=44615.522592592592; Total=44615.522592592592; Msg=$msgOutsideBounds; Wrap=$false }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $workpackage
    $ws.Cells.Item($row, 2).Value = $user
    $ws.Cells.Item($row, 3).Value = $caseId

    $ws.Cells.Item($row, 4).Value = $r.Start
    $ws.Cells.Item($row, 4).NumberFormat = "m/d/yyyy h:mm"
    $ws.Cells.Item($row, 5).Value = $r.End
    $ws.Cells.Item($row, 5).NumberFormat = "m/d/yyyy h:mm"
    $ws.Cells.Item($row, 6).Value = $r.Total
    $ws.Cells.Item($row, 6).NumberFormat = "m/d/yyyy h:mm"

    $ws.Cells.Item($row, 7).Value = $success
    $ws.Cells.Item($row, 8).Value = $r.Msg
    $ws.Cells.Item($row, 9).Value = $notes

    if ($r.Wrap) {
        $ws.Cells.Item($row, 8).WrapText = $true
        $ws.Rows.Item($row).RowHeight = 409.5
    }
}
